$d = $word.ActiveDocument

function Replace-Exact($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2) | Out-Null
}

Replace-Exact "4.49    " "6.78    "
Replace-Exact "2.11    " "5.68    "
Replace-Exact "(5.65)   " "(5.43)   "
Replace-Exact "(5.32)   " "(5.19)   "
Replace-Exact "48.56 ***" "48.64 ***"
Replace-Exact "56.65 ***" "53.11 ***"
Replace-Exact "(11.88)   " "(12.38)   "
Replace-Exact "(11.65)   " "(12.13)   "
Replace-Exact "11.37    " "10.69    "
Replace-Exact "11.04    " "10.44    "
Replace-Exact "(6.14)   " "(6.22)   "
Replace-Exact "(5.62)   " "(5.88)   "
Replace-Exact "7.13    " "7.54    "
Replace-Exact "(6.70)   " "(6.79)   "
Replace-Exact "0.32    " "0.31    "
Replace-Exact "0.47    " "0.42    "
